$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Run Info")

# New data rows appended to the "Run Info" sheet (rows 125-131)
$data = @(
    @(42956.989444444444, "rcp85", 1, 10000, 1, 2.5123978948426076, 0,                  0.3, 0.1, 4, 4, 2, 0.36, 1.5, 0.46, 4.7156000000000002),
    @(42956.989837962959, "rcp85", 1, 10000, 1, 2.4294366844827042, 0,                  0.3, 0.1, 4, 4, 2, 0.36, 1.5, 0.46, 4.7156000000000002),
    @(42956.990289351852, "rcp85", 1, 1000,  2, 2.4243488261286807, 15.384615384615385, 0.3, 0.1, 4, 4, 2, 0.36, 1.5, 0.46, 4.7156000000000002),
    @(42956.990694444445, "rcp85", 1, 1000,  1, 2.318066985015137,  15.384615384615385, 0.3, 0.1, 4, 4, 2, 0.36, 1.5, 0.46, 4.7156000000000002),
    @(42956.99113425926,  "rcp85", 1, 1000,  1, 3.6538281974229836, 15.384615384615385, 0.3, 0.1, 4, 4, 2, 0.36, 1.5, 0.46, 4.7156000000000002),
    @(42956.991215277776, "rcp85", 1, 1000,  1, 3.5865519903925875, 15.384615384615385, 0.3, 0.1, 4, 4, 2, 0.36, 1.5, 0.46, 4.7156000000000002),
    @(42956.991481481484, "rcp85", 1, 1000,  1, 2.5608333463996362, 15.384615384615385, 0.3, 0.1, 4, 4, 2, 0.36, 1.5, 0.46, 4.7156000000000002)
)

$startRow = 125

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]

    # Copy the formatting (number format / style) of the prior row down to the
    # new row before writing its values, so the new cells keep using the
    # existing shared style (e.g. the date style in column A) instead of
    # Excel creating a brand new style entry.
    $ws.Range("A$($row - 1):P$($row - 1)").Copy() | Out-Null
    $ws.Range("A$($row):P$($row)").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
    $excel.CutCopyMode = $false

    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
    $ws.Cells.Item($row, 5).Value = $values[4]
    $ws.Cells.Item($row, 6).Value = $values[5]
    $ws.Cells.Item($row, 7).Value = $values[6]
    $ws.Cells.Item($row, 8).Value = $values[7]
    $ws.Cells.Item($row, 9).Value = $values[8]
    $ws.Cells.Item($row, 10).Value = $values[9]
    $ws.Cells.Item($row, 11).Value = $values[10]
    $ws.Cells.Item($row, 12).Value = $values[11]
    $ws.Cells.Item($row, 13).Value = $values[12]
    $ws.Cells.Item($row, 14).Value = $values[13]
    $ws.Cells.Item($row, 15).Value = $values[14]
    $ws.Cells.Item($row, 16).Value = $values[15]
}

$lastRow = $startRow + $data.Count - 1
$ws.Range("A$($lastRow):P$($lastRow)").Select()
